$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H33").Value = 6565.8125
$ws.Range("I33").Value = 8470.333000000001
$ws.Range("J33").Value = 852.25
$ws.Range("K33").Value = 8470.333000000001
$ws.Range("L33").Value = 852.25
$ws.Range("M33").Value = -8241.333000000001
$ws.Range("N33").Value = -1310.25
$ws.Range("H40").Value = 4626.3335
$ws.Range("I40").Value = 2992.25
$ws.Range("K40").Value = 2992.25
$ws.Range("M40").Value = -2817.25
$ws.Range("H100").Value = 3605.2812
$ws.Range("I100").Value = 3395.5862
$ws.Range("J100").Value = 5632.3335
$ws.Range("K100").Value = 3395.5862
$ws.Range("L100").Value = 5632.3335
$ws.Range("M100").Value = -2854.5862
$ws.Range("N100").Value = -6714.3335
$ws.Range("H112").Value = 4858.4
$ws.Range("J112").Value = 4448
$ws.Range("L112").Value = 13344
$ws.Range("N112").Value = -15560
$ws.Range("H113").Value = 7704.091
$ws.Range("I113").Value = 7733.1665
$ws.Range("J113").Value = 7669.2
$ws.Range("K113").Value = 7733.1665
$ws.Range("L113").Value = 7669.2
$ws.Range("M113").Value = -4479.1665
$ws.Range("N113").Value = -14177.2
$ws.Range("H126").Value = 92299.5
$ws.Range("J126").Value = 92299.5
$ws.Range("L126").Value = 92299.5
$ws.Range("N126").Value = -102179.5
$ws.Range("H132").Value = 6532.2354
$ws.Range("I132").Value = 6969.9
$ws.Range("K132").Value = 20909.7
$ws.Range("M132").Value = -18379.7
$ws.Range("H138").Value = 6980.722
$ws.Range("I138").Value = 7515.9375
$ws.Range("J138").Value = 2699
$ws.Range("K138").Value = 22547.8125
$ws.Range("L138").Value = 8097
$ws.Range("M138").Value = -17407.8125
$ws.Range("N138").Value = -18377

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H61").Value = 5211.25
$ws.Range("J61").Value = 7499.5
$ws.Range("L61").Value = 7499.5
$ws.Range("N61").Value = -7923.5
$ws.Range("H97").Value = 8606.4375
$ws.Range("I97").Value = 12230.777
$ws.Range("J97").Value = 3946.5715
$ws.Range("K97").Value = 12230.777
$ws.Range("L97").Value = 3946.5715
$ws.Range("M97").Value = -11734.777
$ws.Range("N97").Value = -4938.5715
$ws.Range("H125").Value = 40715
$ws.Range("J125").Value = 40715
$ws.Range("L125").Value = 40715
$ws.Range("N125").Value = -50555
$ws.Range("H136").Value = 5211.25
$ws.Range("J136").Value = 7499.5
$ws.Range("L136").Value = 22498.5
$ws.Range("N136").Value = -27598.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H99").Value = 45894.707
$ws.Range("I99").Value = 59235.445
$ws.Range("K99").Value = 59235.445
$ws.Range("M99").Value = -57737.445
$ws.Range("H105").Value = 2777.862
$ws.Range("I105").Value = 2729
$ws.Range("J105").Value = 3437.5
$ws.Range("K105").Value = 2729
$ws.Range("L105").Value = 3437.5
$ws.Range("M105").Value = -982
$ws.Range("N105").Value = -6931.5
$ws.Range("H135").Value = 89995.5
$ws.Range("J135").Value = 89995.5
$ws.Range("L135").Value = 89995.5
$ws.Range("N135").Value = -100135.5
$ws.Range("H138").Value = 98991.8
$ws.Range("J138").Value = 98991.8
$ws.Range("L138").Value = 98991.8
$ws.Range("N138").Value = -109271.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H52").Value = 50000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 50000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 50000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -50588
$ws.Range("H80").Value = 24547.555
$ws.Range("J80").Value = 24547.555
$ws.Range("L80").Value = 24547.555
$ws.Range("N80").Value = -26793.555
$ws.Range("H83").Value = 24547.555
$ws.Range("J83").Value = 24547.555
$ws.Range("L83").Value = 73642.66500000001
$ws.Range("N83").Value = -84874.66500000001
$ws.Range("H99").Value = 4783.1665
$ws.Range("I99").Value = 4766.3335
$ws.Range("J99").Value = 4800
$ws.Range("K99").Value = 4766.3335
$ws.Range("L99").Value = 4800
$ws.Range("M99").Value = -3268.3335
$ws.Range("N99").Value = -7796
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H126").Value = 4783.1665
$ws.Range("I126").Value = 4766.3335
$ws.Range("J126").Value = 4800
$ws.Range("K126").Value = 14299.0005
$ws.Range("L126").Value = 14400
$ws.Range("M126").Value = -11829.0005
$ws.Range("N126").Value = -19340
$ws.Range("H129").Value = 57999.5
$ws.Range("J129").Value = 57999.5
$ws.Range("L129").Value = 57999.5
$ws.Range("N129").Value = -67999.5
$ws.Range("H135").Value = 99779
$ws.Range("J135").Value = 99779
$ws.Range("L135").Value = 99779
$ws.Range("N135").Value = -109919

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H23").Value = 334.7
$ws.Range("J23").Value = 350
$ws.Range("L23").Value = 1050
$ws.Range("N23").Value = -1520
$ws.Range("H122").Value = 464.35715
$ws.Range("J122").Value = 632.7143
$ws.Range("L122").Value = 5694.428699999999
$ws.Range("N122").Value = -10594.4287
$ws.Range("H129").Value = 6995.6523
$ws.Range("I129").Value = 9659.77
$ws.Range("J129").Value = 3532.3
$ws.Range("K129").Value = 28979.31
$ws.Range("L129").Value = 10596.9
$ws.Range("M129").Value = -23979.31
$ws.Range("N129").Value = -20596.9

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H101").Value = 29999
$ws.Range("J101").Value = 29999
$ws.Range("L101").Value = 29999
$ws.Range("N101").Value = -36489
$ws.Range("H126").Value = 6263.5557
$ws.Range("J126").Value = 6878.4
$ws.Range("L126").Value = 20635.2
$ws.Range("N126").Value = -25575.2
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 2181.4
$ws.Range("I122").Value = 1876.75
$ws.Range("K122").Value = 5630.25
$ws.Range("M122").Value = -3180.25
$ws.Range("H129").Value = 69999
$ws.Range("J129").Value = 69999
$ws.Range("L129").Value = 69999
$ws.Range("N129").Value = -79999
